# Auto-generated Excel COM-interop script to apply Golem_Profits sheet updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H57").Value = 45000
$ws.Range("J57").Value = 45000
$ws.Range("L57").Value = 135000
$ws.Range("N57").Value = -135998
$ws.Range("H96").Value = 985.8182
$ws.Range("I96").Value = 942
$ws.Range("J96").Value = 1002.25
$ws.Range("K96").Value = 2826
$ws.Range("L96").Value = 3006.75
$ws.Range("M96").Value = -1453
$ws.Range("N96").Value = -5752.75
$ws.Range("H98").Value = 403.6
$ws.Range("I98").Value = 468.5
$ws.Range("K98").Value = 468.5
$ws.Range("M98").Value = 1029.5
$ws.Range("H107").Value = 53650.234
$ws.Range("I107").Value = 56878.375
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 56878.375
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = -54958.375
$ws.Range("N107").Value = -5840
$ws.Range("H111").Value = 2103.6428
$ws.Range("I111").Value = 2945.25
$ws.Range("K111").Value = 8835.75
$ws.Range("M111").Value = -5768.75
$ws.Range("H112").Value = 1977.7778
$ws.Range("J112").Value = 1977.7778
$ws.Range("L112").Value = 5933.3334
$ws.Range("N112").Value = -8149.3334
$ws.Range("H122").Value = 403.6
$ws.Range("I122").Value = 468.5
$ws.Range("K122").Value = 1405.5
$ws.Range("M122").Value = 1044.5
$ws.Range("H132").Value = 45303.6
$ws.Range("I132").Value = 51578.77
$ws.Range("J132").Value = 4515
$ws.Range("K132").Value = 154736.31
$ws.Range("L132").Value = 13545
$ws.Range("M132").Value = -152206.31
$ws.Range("N132").Value = -18605
$ws.Range("H135").Value = 596.4545
$ws.Range("I135").Value = 575
$ws.Range("J135").Value = 598.6
$ws.Range("K135").Value = 5175
$ws.Range("L135").Value = 5387.400000000001
$ws.Range("M135").Value = -2640
$ws.Range("N135").Value = -10457.4
$ws.Range("H138").Value = 2911.2068
$ws.Range("J138").Value = 3450.739
$ws.Range("L138").Value = 10352.217
$ws.Range("N138").Value = -20632.217
$ws.Range("H141").Value = 250
$ws.Range("I141").Value = 250
$ws.Range("K141").Value = 750
$ws.Range("M141").Value = 4430
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3400
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 3400
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 3400
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -3824
$ws.Range("H110").Value = 22728116
$ws.Range("I110").Value = 785.25
$ws.Range("K110").Value = 785.25
$ws.Range("M110").Value = 1259.75
$ws.Range("H132").Value = 2048
$ws.Range("I132").Value = 1926.2858
$ws.Range("K132").Value = 5778.857400000001
$ws.Range("M132").Value = -3248.857400000001
$ws.Range("H136").Value = 3400
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 3400
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 10200
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -15300
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H22").Value = 724.5
$ws.Range("J22").Value = 275
$ws.Range("L22").Value = 275
$ws.Range("N22").Value = -621
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H81").Value = 22997.5
$ws.Range("J81").Value = 22997.5
$ws.Range("L81").Value = 22997.5
$ws.Range("N81").Value = -25119.5
$ws.Range("H84").Value = 22997.5
$ws.Range("J84").Value = 22997.5
$ws.Range("L84").Value = 68992.5
$ws.Range("N84").Value = -79600.5
$ws.Range("H134").Value = 1880.4286
$ws.Range("I134").Value = 1880.4286
$ws.Range("K134").Value = 5641.2858
$ws.Range("M134").Value = -3106.2858
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 510.83334
$ws.Range("I17").Value = 442
$ws.Range("J17").Value = 560
$ws.Range("K17").Value = 442
$ws.Range("L17").Value = 560
$ws.Range("M17").Value = -268
$ws.Range("N17").Value = -908
$ws.Range("H99").Value = 358591
$ws.Range("I99").Value = 1573.4166
$ws.Range("J99").Value = 2500696.5
$ws.Range("K99").Value = 1573.4166
$ws.Range("L99").Value = 2500696.5
$ws.Range("M99").Value = -75.41660000000002
$ws.Range("N99").Value = -2503692.5
$ws.Range("H100").Value = 100780
$ws.Range("J100").Value = 100780
$ws.Range("L100").Value = 100780
$ws.Range("N100").Value = -102944
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H107").Value = 896.6667
$ws.Range("I107").Value = 821.25
$ws.Range("K107").Value = 821.25
$ws.Range("M107").Value = 1098.75
$ws.Range("H122").Value = 4766
$ws.Range("I122").Value = 1582.625
$ws.Range("J122").Value = 17499.5
$ws.Range("K122").Value = 4747.875
$ws.Range("L122").Value = 52498.5
$ws.Range("M122").Value = -2297.875
$ws.Range("N122").Value = -57398.5
$ws.Range("H126").Value = 358591
$ws.Range("I126").Value = 1573.4166
$ws.Range("J126").Value = 2500696.5
$ws.Range("K126").Value = 4720.2498
$ws.Range("L126").Value = 7502089.5
$ws.Range("M126").Value = -2250.2498
$ws.Range("N126").Value = -7507029.5
$ws.Range("H132").Value = 754.8571
$ws.Range("I132").Value = 547.3333
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 1641.9999
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = 888.0001
$ws.Range("N132").Value = -11060
$ws.Range("H141").Value = 429696.2
$ws.Range("J141").Value = 761109.7
$ws.Range("L141").Value = 761109.7
$ws.Range("N141").Value = -771469.7
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 631.8571
$ws.Range("J113").Value = 525.6667
$ws.Range("L113").Value = 1577.0001
$ws.Range("N113").Value = -5917.0001
$ws.Range("H128").Value = 629997
$ws.Range("I128").Value = 629997
$ws.Range("K128").Value = 1889991
$ws.Range("M128").Value = -1885011
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 17253.572
$ws.Range("I43").Value = 3943.75
$ws.Range("K43").Value = 3943.75
$ws.Range("M43").Value = -3792.75
$ws.Range("H46").Value = 12799.8
$ws.Range("I46").Value = 12799.8
$ws.Range("K46").Value = 12799.8
$ws.Range("M46").Value = -12643.8
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H102").Value = 1282.125
$ws.Range("I102").Value = 1036.7142
$ws.Range("K102").Value = 1036.7142
$ws.Range("M102").Value = 585.2858000000001
$ws.Range("H127").Value = 30000
$ws.Range("I127").Value = 30000
$ws.Range("K127").Value = 30000
$ws.Range("M127").Value = -25040
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7999.6
$ws.Range("I7").Value = 7999.6
$ws.Range("K7").Value = 7999.6
$ws.Range("M7").Value = -7887.6
$ws.Range("H16").Value = 1000
$ws.Range("I16").Value = 1000
$ws.Range("K16").Value = 1000
$ws.Range("M16").Value = -830
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("N27").ClearContents()
$ws.Range("H40").Value = 2286
$ws.Range("I40").Value = 2286
$ws.Range("K40").Value = 2286
$ws.Range("M40").Value = -2150
$ws.Range("H126").Value = 7999.6
$ws.Range("I126").Value = 7999.6
$ws.Range("K126").Value = 23998.8
$ws.Range("M126").Value = -21528.8
$ws.Range("H132").Value = 1765.5
$ws.Range("I132").Value = 975.4286
$ws.Range("J132").Value = 2555.5715
$ws.Range("K132").Value = 2926.2858
$ws.Range("L132").Value = 7666.7145
$ws.Range("M132").Value = -396.2857999999997
$ws.Range("N132").Value = -12726.7145
$ws.Range("H136").Value = 599.5
$ws.Range("I136").Value = 599.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 1798.5
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = 751.5
$ws.Range("N136").ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1605.0834
$ws.Range("I96").Value = 1172.8334
$ws.Range("J96").Value = 2037.3334
$ws.Range("K96").Value = 1172.8334
$ws.Range("L96").Value = 2037.3334
$ws.Range("M96").Value = 200.1666
$ws.Range("N96").Value = -4783.3334
$ws.Range("H122").Value = 1185.8182
$ws.Range("I122").Value = 927.1111
$ws.Range("K122").Value = 2781.3333
$ws.Range("M122").Value = -331.3332999999998
$ws.Range("H126").Value = 3484.5
$ws.Range("I126").Value = 2760
$ws.Range("K126").Value = 8280
$ws.Range("M126").Value = -5810
$ws.Range("H132").Value = 1355.25
$ws.Range("I132").Value = 976.7143
$ws.Range("J132").Value = 4005
$ws.Range("K132").Value = 2930.1429
$ws.Range("L132").Value = 12015
$ws.Range("M132").Value = -400.1428999999998
$ws.Range("N132").Value = -17075
$ws.Range("H136").Value = 3310.8
$ws.Range("I136").Value = 2183
$ws.Range("K136").Value = 6549
$ws.Range("M136").Value = -3999
